# Update cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text value.
# Values that parse as plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (quote-prefixed) instead of
# converting them to numbers (which would lose the original formatting,
# e.g. "1.00" becoming 1, or introduce floating point artifacts).
$updates = @(
    @("D2", "58.601.70"),
    @("E2", "  -4.91%  "),
    @("D3", "2.470.06"),
    @("E3", "  -3.78%  "),
    @("E4", "  -0.09%  "),
    @("D5", "'529.75"),
    @("E5", "  -3.64%  "),
    @("D6", "'142.10"),
    @("E6", "  -7.86%  "),
    @("D7", "'0.997"),
    @("E7", "  -0.35%  "),
    @("D8", "'0.567"),
    @("E8", "  -5.13%  "),
    @("D9", "2.490.19"),
    @("E9", "  -2.94%  "),
    @("D10", "'0.0987"),
    @("E10", "  -5.18%  "),
    @("E11", "  -2.80%  "),
    @("D12", "'5.44"),
    @("E12", "  +0.44%  "),
    @("D13", "'0.347"),
    @("E13", "  -4.71%  "),
    @("D14", "2.900.17"),
    @("E14", "  -4.04%  "),
    @("D15", "'23.49"),
    @("E15", "  -7.57%  "),
    @("D16", "58.453.14"),
    @("E16", "  -5.00%  "),
    @("D17", "'0.0000136"),
    @("E17", "  -4.99%  "),
    @("D18", "2.467.25"),
    @("E18", "  -4.18%  "),
    @("D19", "'11.20"),
    @("E19", "  -2.99%  "),
    @("D20", "'4.23"),
    @("E20", "  -6.63%  "),
    @("D21", "'319.86"),
    @("E21", "  -5.14%  "),
    @("E22", "  -0.26%  "),
    @("D23", "'5.68"),
    @("E23", "  -6.07%  "),
    @("D24", "'60.19"),
    @("E24", "  -4.68%  "),
    @("D25", "'0.429"),
    @("E25", "  -13.06%  "),
    @("D26", "'1.00"),
    @("E26", "  +0.39%  "),
    @("B27", "Kaspa"),
    @("C27", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"),
    @("D27", "'0.159"),
    @("E27", "  -5.25%  "),
    @("B28", "WrappedeETH"),
    @("C28", "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"),
    @("D28", "2.572.16"),
    @("E28", "  -4.62%  "),
    @("D29", "'7.65"),
    @("E29", "  -5.56%  "),
    @("D30", "'6.68"),
    @("E30", "  -9.49%  "),
    @("D31", "'1.76"),
    @("E31", "  -7.77%  "),
    @("D32", "0.0₃0751"),
    @("E32", "  -10.15%  "),
    @("D33", "'1.21"),
    @("E33", "  -6.03%  "),
    @("D34", "'0.996"),
    @("E34", "  -0.32%  "),
    @("D35", "'156.66"),
    @("E35", "  -1.62%  "),
    @("D36", "'1.39"),
    @("E36", "  -0.69%  "),
    @("D37", "'18.40"),
    @("E37", "  -3.99%  "),
    @("D38", "'4.28"),
    @("E38", "  -9.04%  "),
    @("D39", "'1.57"),
    @("E39", "  -13.04%  "),
    @("D40", "'5.70"),
    @("E40", "  -4.46%  "),
    @("D41", "'36.57"),
    @("E41", "  -2.45%  "),
    @("D42", "'300.03"),
    @("E42", "  -10.07%  "),
    @("D43", "'3.60"),
    @("E43", "  -8.88%  "),
    @("B44", "SuiNetwork"),
    @("C44", "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"),
    @("D44", "'0.773"),
    @("E44", "  -17.75%  "),
    @("B45", "FirstDigitalUSD"),
    @("C45", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"),
    @("D45", "'0.995"),
    @("E45", "  -0.30%  "),
    @("D46", "'10.74"),
    @("E46", "  -1.73%  "),
    @("D47", "'0.590"),
    @("E47", "  -2.21%  "),
    @("D48", "'123.78"),
    @("E48", "  +0.56%  "),
    @("D49", "'0.0919"),
    @("E49", "  -4.75%  "),
    @("D50", "'18.38"),
    @("E50", "  -5.79%  "),
    @("D51", "'0.0510"),
    @("E51", "  -6.57%  ")
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
